$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.120.81"
$ws.Range("D3").Value = "1.654.94"
$ws.Range("E3").Value = "  -0.20%  "
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5250"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.85%  "
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2609"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.67%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06352"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.38"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07795"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.11%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.683.70"
$ws.Range("E12").Value = "  +1.62%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.504"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.95%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5479"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.08%  "
$ws.Range("D15").Value = "0.0₅8202"
$ws.Range("E15").Value = "  +0.99%  "
$ws.Range("E16").Value = "  +1.59%  "
$ws.Range("D17").Value = "26.116.65"
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.579"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "191.73"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("E21").Value = "  +0.14%  "
$ws.Range("E22").Value = "  -0.44%  "
$ws.Range("E23").Value = "  -0.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "142.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.78%  "
$ws.Range("E25").Value = "  +1.54%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.260"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.21%  "
$ws.Range("E27").Value = "  +0.67%  "
$ws.Range("E28").Value = "  +1.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05918"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.279"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.67%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.526"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.59%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.252"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.590"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9535"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.786"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.62%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.410"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.29%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5697"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.53%  "
$ws.Range("E38").Value = "  +1.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.804"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8500"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.58%  "
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "1.031.76"
$ws.Range("E42").Value = "  +1.95%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "103.15"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.97%  "
$ws.Range("D44").Value = "1.800.15"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "57.23"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.87%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.008"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4301"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.42%  "
$ws.Range("E48").Value = "  +1.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05166"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.862"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.09698"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.13%  "
